$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at row 5 to make room for the new FAPs sending-cluster block
# (existing "sCs" rows 5-7 shift down to rows 8-10)
$ws.Rows("5:7").Insert()

# Row 2
$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "Dll1"
$ws.Range("C2").Value2 = "Notch2"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 5.190862666666667
$ws.Range("H2").Value2 = 15.572588
$ws.Range("I2").Value2 = 0.8740249884703439
$ws.Range("J2").Value2 = 0.874024988470344
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 33.380049
$ws.Range("N2").Value2 = 100.140147
$ws.Range("O2").Value2 = 0.3891462059670435
$ws.Range("P2").Value2 = 0.3891462059670435
$ws.Range("Q2").Value2 = 173.271250165604
$ws.Range("R2").Value2 = 1559.441251490436
$ws.Range("S2").Value2 = 0.3401235081836233
$ws.Range("T2").Value2 = 0.3401235081836234

# Row 3
$ws.Range("A3").Value2 = "ECs"
$ws.Range("B3").Value2 = "Dll1"
$ws.Range("C3").Value2 = "Notch2"
$ws.Range("D3").Value2 = "FAPs"
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 5.190862666666667
$ws.Range("H3").Value2 = 15.572588
$ws.Range("I3").Value2 = 0.8740249884703439
$ws.Range("J3").Value2 = 0.874024988470344
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 33.85786133333334
$ws.Range("N3").Value2 = 101.573584
$ws.Range("O3").Value2 = 0.3947165649764305
$ws.Range("P3").Value2 = 0.3947165649764305
$ws.Range("Q3").Value2 = 175.7515083683769
$ws.Range("R3").Value2 = 1581.763575315392
$ws.Range("S3").Value2 = 0.3449921411525785
$ws.Range("T3").Value2 = 0.3449921411525785

# Row 4
$ws.Range("A4").Value2 = "ECs"
$ws.Range("B4").Value2 = "Dll1"
$ws.Range("C4").Value2 = "Notch2"
$ws.Range("D4").Value2 = "sCs"
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 5.190862666666667
$ws.Range("H4").Value2 = 15.572588
$ws.Range("I4").Value2 = 0.8740249884703439
$ws.Range("J4").Value2 = 0.874024988470344
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 18.53974466666667
$ws.Range("N4").Value2 = 55.61923400000001
$ws.Range("O4").Value2 = 0.2161372290565261
$ws.Range("P4").Value2 = 0.2161372290565261
$ws.Range("Q4").Value2 = 96.23726843973245
$ws.Range("R4").Value2 = 866.1354159575922
$ws.Range("S4").Value2 = 0.1889093391341423
$ws.Range("T4").Value2 = 0.1889093391341423

# Row 5
$ws.Range("A5").Value2 = "FAPs"
$ws.Range("B5").Value2 = "Dll1"
$ws.Range("C5").Value2 = "Notch2"
$ws.Range("D5").Value2 = "ECs"
$ws.Range("E5").Value2 = 2
$ws.Range("F5").Value2 = 0.6666666666666666
$ws.Range("G5").Value2 = 0.100996
$ws.Range("H5").Value2 = 0.302988
$ws.Range("I5").Value2 = 0.0170054639091879
$ws.Range("J5").Value2 = 0.0170054639091879
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 33.380049
$ws.Range("N5").Value2 = 100.140147
$ws.Range("O5").Value2 = 0.3891462059670435
$ws.Range("P5").Value2 = 0.3891462059670435
$ws.Range("Q5").Value2 = 3.371251428803999
$ws.Range("R5").Value2 = 30.341262859236
$ws.Range("S5").Value2 = 0.006617611760969959
$ws.Range("T5").Value2 = 0.006617611760969959

# Row 6
$ws.Range("A6").Value2 = "FAPs"
$ws.Range("B6").Value2 = "Dll1"
$ws.Range("C6").Value2 = "Notch2"
$ws.Range("D6").Value2 = "FAPs"
$ws.Range("E6").Value2 = 2
$ws.Range("F6").Value2 = 0.6666666666666666
$ws.Range("G6").Value2 = 0.100996
$ws.Range("H6").Value2 = 0.302988
$ws.Range("I6").Value2 = 0.0170054639091879
$ws.Range("J6").Value2 = 0.0170054639091879
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 33.85786133333334
$ws.Range("N6").Value2 = 101.573584
$ws.Range("O6").Value2 = 0.3947165649764305
$ws.Range("P6").Value2 = 0.3947165649764305
$ws.Range("Q6").Value2 = 3.419508563221334
$ws.Range("R6").Value2 = 30.775577068992
$ws.Range("S6").Value2 = 0.006712338300065309
$ws.Range("T6").Value2 = 0.006712338300065309

# Row 7
$ws.Range("A7").Value2 = "FAPs"
$ws.Range("B7").Value2 = "Dll1"
$ws.Range("C7").Value2 = "Notch2"
$ws.Range("D7").Value2 = "sCs"
$ws.Range("E7").Value2 = 2
$ws.Range("F7").Value2 = 0.6666666666666666
$ws.Range("G7").Value2 = 0.100996
$ws.Range("H7").Value2 = 0.302988
$ws.Range("I7").Value2 = 0.0170054639091879
$ws.Range("J7").Value2 = 0.0170054639091879
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 18.53974466666667
$ws.Range("N7").Value2 = 55.61923400000001
$ws.Range("O7").Value2 = 0.2161372290565261
$ws.Range("P7").Value2 = 0.2161372290565261
$ws.Range("Q7").Value2 = 1.872440052354666
$ws.Range("R7").Value2 = 16.851960471192
$ws.Range("S7").Value2 = 0.003675513848152632
$ws.Range("T7").Value2 = 0.003675513848152632

# Row 8
$ws.Range("A8").Value2 = "sCs"
$ws.Range("B8").Value2 = "Dll1"
$ws.Range("C8").Value2 = "Notch2"
$ws.Range("D8").Value2 = "ECs"
$ws.Range("E8").Value2 = 3
$ws.Range("F8").Value2 = 1
$ws.Range("G8").Value2 = 0.6471736666666666
$ws.Range("H8").Value2 = 1.941521
$ws.Range("I8").Value2 = 0.1089695476204681
$ws.Range("J8").Value2 = 0.1089695476204681
$ws.Range("K8").Value2 = 3
$ws.Range("L8").Value2 = 1
$ws.Range("M8").Value2 = 33.380049
$ws.Range("N8").Value2 = 100.140147
$ws.Range("O8").Value2 = 0.3891462059670435
$ws.Range("P8").Value2 = 0.3891462059670435
$ws.Range("Q8").Value2 = 21.602688704843
$ws.Range("R8").Value2 = 194.424198343587
$ws.Range("S8").Value2 = 0.04240508602245025
$ws.Range("T8").Value2 = 0.04240508602245025

# Row 9
$ws.Range("A9").Value2 = "sCs"
$ws.Range("B9").Value2 = "Dll1"
$ws.Range("C9").Value2 = "Notch2"
$ws.Range("D9").Value2 = "FAPs"
$ws.Range("E9").Value2 = 3
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = 0.6471736666666666
$ws.Range("H9").Value2 = 1.941521
$ws.Range("I9").Value2 = 0.1089695476204681
$ws.Range("J9").Value2 = 0.1089695476204681
$ws.Range("K9").Value2 = 3
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 33.85786133333334
$ws.Range("N9").Value2 = 101.573584
$ws.Range("O9").Value2 = 0.3947165649764305
$ws.Range("P9").Value2 = 0.3947165649764305
$ws.Range("Q9").Value2 = 21.91191626458489
$ws.Range("R9").Value2 = 197.207246381264
$ws.Range("S9").Value2 = 0.04301208552378676
$ws.Range("T9").Value2 = 0.04301208552378675

# Row 10
$ws.Range("A10").Value2 = "sCs"
$ws.Range("B10").Value2 = "Dll1"
$ws.Range("C10").Value2 = "Notch2"
$ws.Range("D10").Value2 = "sCs"
$ws.Range("E10").Value2 = 3
$ws.Range("F10").Value2 = 1
$ws.Range("G10").Value2 = 0.6471736666666666
$ws.Range("H10").Value2 = 1.941521
$ws.Range("I10").Value2 = 0.1089695476204681
$ws.Range("J10").Value2 = 0.1089695476204681
$ws.Range("K10").Value2 = 3
$ws.Range("L10").Value2 = 1
$ws.Range("M10").Value2 = 18.53974466666667
$ws.Range("N10").Value2 = 55.61923400000001
$ws.Range("O10").Value2 = 0.2161372290565261
$ws.Range("P10").Value2 = 0.2161372290565261
$ws.Range("Q10").Value2 = 11.99843453499044
$ws.Range("R10").Value2 = 107.985910814914
$ws.Range("S10").Value2 = 0.02355237607423115
$ws.Range("T10").Value2 = 0.02355237607423115
